$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A 68622-2021
$ws.Range("A2").Value = 'A 68622-2021'
$ws.Range("B2").Value2 = 44529.0
$ws.Range("C2").Value2 = 46060.0
$ws.Range("D2").Value = 'VÄRMLANDS LÄN'
$ws.Range("E2").Value = 'HAMMARÖ'
$ws.Range("G2").Value2 = 9.7
$ws.Range("H2").Value2 = 5
$ws.Range("I2").Value2 = 1
$ws.Range("J2").Value2 = 1
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0
$ws.Range("M2").Value2 = 0
$ws.Range("N2").Value2 = 0
$ws.Range("O2").Value2 = 2
$ws.Range("P2").Value2 = 1
$ws.Range("Q2").Value2 = 7
$ws.Range("R2").Value = 'Knärot`r`nTallticka`r`nJättesvampmal`r`nÅkergroda`r`nKopparödla`r`nVanlig groda`r`nVanlig padda'
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 68622-2021 artfynd.xlsx", "A 68622-2021")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 68622-2021 karta.png", "A 68622-2021")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/knärot/A 68622-2021 karta knärot.png", "A 68622-2021")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 68622-2021 FSC-klagomål.docx", "A 68622-2021")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 68622-2021 FSC-klagomål mail.docx", "A 68622-2021")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 68622-2021 tillsynsbegäran.docx", "A 68622-2021")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 68622-2021 tillsynsbegäran mail.docx", "A 68622-2021")'

# Row 3: A 61380-2023
$ws.Range("A3").Value = 'A 61380-2023'
$ws.Range("B3").Value2 = 45264.0
$ws.Range("C3").Value2 = 46060.0
$ws.Range("D3").Value = 'VÄRMLANDS LÄN'
$ws.Range("E3").Value = 'HAMMARÖ'
$ws.Range("G3").Value2 = 3.3
$ws.Range("H3").Value2 = 1
$ws.Range("I3").Value2 = 1
$ws.Range("J3").Value2 = 3
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0
$ws.Range("M3").Value2 = 0
$ws.Range("N3").Value2 = 0
$ws.Range("O3").Value2 = 4
$ws.Range("P3").Value2 = 1
$ws.Range("Q3").Value2 = 5
$ws.Range("R3").Value = 'Knärot`r`nBlå taggsvamp`r`nGropticka`r`nMotaggsvamp`r`nGullgröppa'
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 61380-2023 artfynd.xlsx", "A 61380-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 61380-2023 karta.png", "A 61380-2023")'
$ws.Range("U3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/knärot/A 61380-2023 karta knärot.png", "A 61380-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 61380-2023 FSC-klagomål.docx", "A 61380-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 61380-2023 FSC-klagomål mail.docx", "A 61380-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 61380-2023 tillsynsbegäran.docx", "A 61380-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 61380-2023 tillsynsbegäran mail.docx", "A 61380-2023")'
$ws.Range("Z3").ClearContents()

# Row 4: A 47571-2025
$ws.Range("A4").Value = 'A 47571-2025'
$ws.Range("B4").Value2 = 45931.0
$ws.Range("C4").Value2 = 46060.0
$ws.Range("D4").Value = 'VÄRMLANDS LÄN'
$ws.Range("E4").Value = 'HAMMARÖ'
$ws.Range("G4").Value2 = 8.3
$ws.Range("H4").Value2 = 4
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = 1
$ws.Range("K4").Value2 = 2
$ws.Range("L4").Value2 = 0
$ws.Range("M4").Value2 = 0
$ws.Range("N4").Value2 = 0
$ws.Range("O4").Value2 = 3
$ws.Range("P4").Value2 = 2
$ws.Range("Q4").Value2 = 5
$ws.Range("R4").Value = 'Knärot`r`nLångskägg`r`nMotaggsvamp`r`nKungsfågel`r`nTjäder'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 47571-2025 artfynd.xlsx", "A 47571-2025")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 47571-2025 karta.png", "A 47571-2025")'
$ws.Range("U4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/knärot/A 47571-2025 karta knärot.png", "A 47571-2025")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 47571-2025 FSC-klagomål.docx", "A 47571-2025")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 47571-2025 FSC-klagomål mail.docx", "A 47571-2025")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 47571-2025 tillsynsbegäran.docx", "A 47571-2025")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 47571-2025 tillsynsbegäran mail.docx", "A 47571-2025")'
$ws.Range("Z4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/fåglar/A 47571-2025 prioriterade fågelarter.docx", "A 47571-2025")'

# Row 5: A 24616-2022
$ws.Range("A5").Value = 'A 24616-2022'
$ws.Range("B5").Value2 = 44727.0
$ws.Range("C5").Value2 = 46060.0
$ws.Range("D5").Value = 'VÄRMLANDS LÄN'
$ws.Range("E5").Value = 'HAMMARÖ'
$ws.Range("F5").Value = 'Kommuner'
$ws.Range("G5").Value2 = 4
$ws.Range("H5").Value2 = 3
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = 2
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0
$ws.Range("M5").Value2 = 0
$ws.Range("N5").Value2 = 0
$ws.Range("O5").Value2 = 3
$ws.Range("P5").Value2 = 1
$ws.Range("Q5").Value2 = 3
$ws.Range("R5").Value = 'Knärot`r`nSpillkråka`r`nTretåig hackspett'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 24616-2022 artfynd.xlsx", "A 24616-2022")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 24616-2022 karta.png", "A 24616-2022")'
$ws.Range("U5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/knärot/A 24616-2022 karta knärot.png", "A 24616-2022")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 24616-2022 FSC-klagomål.docx", "A 24616-2022")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 24616-2022 FSC-klagomål mail.docx", "A 24616-2022")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 24616-2022 tillsynsbegäran.docx", "A 24616-2022")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 24616-2022 tillsynsbegäran mail.docx", "A 24616-2022")'
$ws.Range("Z5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/fåglar/A 24616-2022 prioriterade fågelarter.docx", "A 24616-2022")'

# Row 6: A 3573-2026
$ws.Range("A6").Value = 'A 3573-2026'
$ws.Range("B6").Value2 = 46042.69962962963
$ws.Range("C6").Value2 = 46060.0
$ws.Range("D6").Value = 'VÄRMLANDS LÄN'
$ws.Range("E6").Value = 'HAMMARÖ'
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value2 = 4.3
$ws.Range("H6").Value2 = 3
$ws.Range("I6").Value2 = 0
$ws.Range("J6").Value2 = 0
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0
$ws.Range("M6").Value2 = 0
$ws.Range("N6").Value2 = 0
$ws.Range("O6").Value2 = 1
$ws.Range("P6").Value2 = 1
$ws.Range("Q6").Value2 = 3
$ws.Range("R6").Value = 'Knärot`r`nVanlig groda`r`nVanlig padda'
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 3573-2026 artfynd.xlsx", "A 3573-2026")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 3573-2026 karta.png", "A 3573-2026")'
$ws.Range("U6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/knärot/A 3573-2026 karta knärot.png", "A 3573-2026")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 3573-2026 FSC-klagomål.docx", "A 3573-2026")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 3573-2026 FSC-klagomål mail.docx", "A 3573-2026")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 3573-2026 tillsynsbegäran.docx", "A 3573-2026")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 3573-2026 tillsynsbegäran mail.docx", "A 3573-2026")'
$ws.Range("Z6").ClearContents()

# Row 7: A 53276-2023
$ws.Range("A7").Value = 'A 53276-2023'
$ws.Range("B7").Value2 = 45229.0
$ws.Range("C7").Value2 = 46060.0
$ws.Range("D7").Value = 'VÄRMLANDS LÄN'
$ws.Range("E7").Value = 'HAMMARÖ'
$ws.Range("F7").Value = 'Övriga Aktiebolag'
$ws.Range("G7").Value2 = 0.3
$ws.Range("H7").Value2 = 2
$ws.Range("I7").Value2 = 0
$ws.Range("J7").Value2 = 0
$ws.Range("K7").Value2 = 0
$ws.Range("L7").Value2 = 0
$ws.Range("M7").Value2 = 0
$ws.Range("N7").Value2 = 0
$ws.Range("O7").Value2 = 0
$ws.Range("P7").Value2 = 0
$ws.Range("Q7").Value2 = 2
$ws.Range("R7").Value = 'Mindre flugsnappare`r`nVanlig groda'
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 53276-2023 artfynd.xlsx", "A 53276-2023")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 53276-2023 karta.png", "A 53276-2023")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 53276-2023 FSC-klagomål.docx", "A 53276-2023")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 53276-2023 FSC-klagomål mail.docx", "A 53276-2023")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 53276-2023 tillsynsbegäran.docx", "A 53276-2023")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 53276-2023 tillsynsbegäran mail.docx", "A 53276-2023")'
$ws.Range("Z7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/fåglar/A 53276-2023 prioriterade fågelarter.docx", "A 53276-2023")'

# Row 8: A 20755-2021
$ws.Range("A8").Value = 'A 20755-2021'
$ws.Range("B8").Value2 = 44316.0
$ws.Range("C8").Value2 = 46060.0
$ws.Range("D8").Value = 'VÄRMLANDS LÄN'
$ws.Range("E8").Value = 'HAMMARÖ'
$ws.Range("G8").Value2 = 16.8
$ws.Range("H8").Value2 = 1
$ws.Range("I8").Value2 = 0
$ws.Range("J8").Value2 = 1
$ws.Range("K8").Value2 = 0
$ws.Range("L8").Value2 = 0
$ws.Range("M8").Value2 = 0
$ws.Range("N8").Value2 = 0
$ws.Range("O8").Value2 = 1
$ws.Range("P8").Value2 = 0
$ws.Range("Q8").Value2 = 1
$ws.Range("R8").Value = 'Gulsparv'
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 20755-2021 artfynd.xlsx", "A 20755-2021")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 20755-2021 karta.png", "A 20755-2021")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 20755-2021 FSC-klagomål.docx", "A 20755-2021")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 20755-2021 FSC-klagomål mail.docx", "A 20755-2021")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 20755-2021 tillsynsbegäran.docx", "A 20755-2021")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 20755-2021 tillsynsbegäran mail.docx", "A 20755-2021")'
$ws.Range("Z8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/fåglar/A 20755-2021 prioriterade fågelarter.docx", "A 20755-2021")'

# Row 9: A 58382-2024
$ws.Range("A9").Value = 'A 58382-2024'
$ws.Range("B9").Value2 = 45632.0
$ws.Range("C9").Value2 = 46060.0
$ws.Range("D9").Value = 'VÄRMLANDS LÄN'
$ws.Range("E9").Value = 'HAMMARÖ'
$ws.Range("G9").Value2 = 1.6
$ws.Range("H9").Value2 = 1
$ws.Range("I9").Value2 = 0
$ws.Range("J9").Value2 = 0
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0
$ws.Range("M9").Value2 = 0
$ws.Range("N9").Value2 = 0
$ws.Range("O9").Value2 = 1
$ws.Range("P9").Value2 = 1
$ws.Range("Q9").Value2 = 1
$ws.Range("R9").Value = 'Knärot'
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 58382-2024 artfynd.xlsx", "A 58382-2024")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 58382-2024 karta.png", "A 58382-2024")'
$ws.Range("U9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/knärot/A 58382-2024 karta knärot.png", "A 58382-2024")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 58382-2024 FSC-klagomål.docx", "A 58382-2024")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 58382-2024 FSC-klagomål mail.docx", "A 58382-2024")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 58382-2024 tillsynsbegäran.docx", "A 58382-2024")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 58382-2024 tillsynsbegäran mail.docx", "A 58382-2024")'

# Row 10: A 45407-2025
$ws.Range("A10").Value = 'A 45407-2025'
$ws.Range("B10").Value2 = 45922.0
$ws.Range("C10").Value2 = 46060.0
$ws.Range("D10").Value = 'VÄRMLANDS LÄN'
$ws.Range("E10").Value = 'HAMMARÖ'
$ws.Range("G10").Value2 = 11.6
$ws.Range("H10").Value2 = 1
$ws.Range("I10").Value2 = 0
$ws.Range("J10").Value2 = 0
$ws.Range("K10").Value2 = 0
$ws.Range("L10").Value2 = 0
$ws.Range("M10").Value2 = 0
$ws.Range("N10").Value2 = 0
$ws.Range("O10").Value2 = 0
$ws.Range("P10").Value2 = 0
$ws.Range("Q10").Value2 = 1
$ws.Range("R10").Value = 'Vanlig padda'
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 45407-2025 artfynd.xlsx", "A 45407-2025")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 45407-2025 karta.png", "A 45407-2025")'
$ws.Range("U10").ClearContents()
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 45407-2025 FSC-klagomål.docx", "A 45407-2025")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 45407-2025 FSC-klagomål mail.docx", "A 45407-2025")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 45407-2025 tillsynsbegäran.docx", "A 45407-2025")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 45407-2025 tillsynsbegäran mail.docx", "A 45407-2025")'

# Row 11: A 24618-2022
$ws.Range("A11").Value = 'A 24618-2022'
$ws.Range("B11").Value2 = 44727.0
$ws.Range("C11").Value2 = 46060.0
$ws.Range("D11").Value = 'VÄRMLANDS LÄN'
$ws.Range("E11").Value = 'HAMMARÖ'
$ws.Range("F11").Value = 'Kommuner'
$ws.Range("G11").Value2 = 4.6
$ws.Range("H11").Value2 = 1
$ws.Range("I11").Value2 = 0
$ws.Range("J11").Value2 = 0
$ws.Range("K11").Value2 = 1
$ws.Range("L11").Value2 = 0
$ws.Range("M11").Value2 = 0
$ws.Range("N11").Value2 = 0
$ws.Range("O11").Value2 = 1
$ws.Range("P11").Value2 = 1
$ws.Range("Q11").Value2 = 1
$ws.Range("R11").Value = 'Knärot'
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 24618-2022 artfynd.xlsx", "A 24618-2022")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 24618-2022 karta.png", "A 24618-2022")'
$ws.Range("U11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/knärot/A 24618-2022 karta knärot.png", "A 24618-2022")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 24618-2022 FSC-klagomål.docx", "A 24618-2022")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 24618-2022 FSC-klagomål mail.docx", "A 24618-2022")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 24618-2022 tillsynsbegäran.docx", "A 24618-2022")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 24618-2022 tillsynsbegäran mail.docx", "A 24618-2022")'

# Row 12: A 7593-2025
$ws.Range("A12").Value = 'A 7593-2025'
$ws.Range("B12").Value2 = 45705.0
$ws.Range("C12").Value2 = 46060.0
$ws.Range("D12").Value = 'VÄRMLANDS LÄN'
$ws.Range("E12").Value = 'HAMMARÖ'
$ws.Range("G12").Value2 = 1.3
$ws.Range("H12").Value2 = 1
$ws.Range("I12").Value2 = 0
$ws.Range("J12").Value2 = 0
$ws.Range("K12").Value2 = 0
$ws.Range("L12").Value2 = 0
$ws.Range("M12").Value2 = 0
$ws.Range("N12").Value2 = 0
$ws.Range("O12").Value2 = 0
$ws.Range("P12").Value2 = 0
$ws.Range("Q12").Value2 = 1
$ws.Range("R12").Value = 'Trollpipistrell'
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 7593-2025 artfynd.xlsx", "A 7593-2025")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 7593-2025 karta.png", "A 7593-2025")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 7593-2025 FSC-klagomål.docx", "A 7593-2025")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 7593-2025 FSC-klagomål mail.docx", "A 7593-2025")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 7593-2025 tillsynsbegäran.docx", "A 7593-2025")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 7593-2025 tillsynsbegäran mail.docx", "A 7593-2025")'
$ws.Range("Z12").ClearContents()

# Row 13: A 55068-2023
$ws.Range("A13").Value = 'A 55068-2023'
$ws.Range("B13").Value2 = 45237.0
$ws.Range("C13").Value2 = 46060.0
$ws.Range("D13").Value = 'VÄRMLANDS LÄN'
$ws.Range("E13").Value = 'HAMMARÖ'
$ws.Range("F13").Value = 'Kommuner'
$ws.Range("G13").Value2 = 2.4
$ws.Range("H13").Value2 = 0
$ws.Range("I13").Value2 = 0
$ws.Range("J13").Value2 = 1
$ws.Range("K13").Value2 = 0
$ws.Range("L13").Value2 = 0
$ws.Range("M13").Value2 = 0
$ws.Range("N13").Value2 = 0
$ws.Range("O13").Value2 = 1
$ws.Range("P13").Value2 = 0
$ws.Range("Q13").Value2 = 1
$ws.Range("R13").Value = 'Motaggsvamp'
$ws.Range("S13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/artfynd/A 55068-2023 artfynd.xlsx", "A 55068-2023")'
$ws.Range("T13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/kartor/A 55068-2023 karta.png", "A 55068-2023")'
$ws.Range("U13").ClearContents()
$ws.Range("V13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomål/A 55068-2023 FSC-klagomål.docx", "A 55068-2023")'
$ws.Range("W13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/klagomålsmail/A 55068-2023 FSC-klagomål mail.docx", "A 55068-2023")'
$ws.Range("X13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsyn/A 55068-2023 tillsynsbegäran.docx", "A 55068-2023")'
$ws.Range("Y13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1761/tillsynsmail/A 55068-2023 tillsynsbegäran mail.docx", "A 55068-2023")'

# Row 14: A 16890-2021
$ws.Range("A14").Value = 'A 16890-2021'
$ws.Range("B14").Value2 = 44295.0
$ws.Range("C14").Value2 = 46060.0
$ws.Range("D14").Value = 'VÄRMLANDS LÄN'
$ws.Range("E14").Value = 'HAMMARÖ'
$ws.Range("G14").Value2 = 1.7
$ws.Range("H14").Value2 = 0
$ws.Range("I14").Value2 = 0
$ws.Range("J14").Value2 = 0
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = 0
$ws.Range("M14").Value2 = 0
$ws.Range("N14").Value2 = 0
$ws.Range("O14").Value2 = 0
$ws.Range("P14").Value2 = 0
$ws.Range("Q14").Value2 = 0

# Row 15: A 56835-2021
$ws.Range("A15").Value = 'A 56835-2021'
$ws.Range("B15").Value2 = 44481.0
$ws.Range("C15").Value2 = 46060.0
$ws.Range("D15").Value = 'VÄRMLANDS LÄN'
$ws.Range("E15").Value = 'HAMMARÖ'
$ws.Range("F15").Value = 'Kommuner'
$ws.Range("G15").Value2 = 18.8
$ws.Range("H15").Value2 = 0
$ws.Range("I15").Value2 = 0
$ws.Range("J15").Value2 = 0
$ws.Range("K15").Value2 = 0
$ws.Range("L15").Value2 = 0
$ws.Range("M15").Value2 = 0
$ws.Range("N15").Value2 = 0
$ws.Range("O15").Value2 = 0
$ws.Range("P15").Value2 = 0
$ws.Range("Q15").Value2 = 0

# Row 16: A 42957-2023
$ws.Range("A16").Value = 'A 42957-2023'
$ws.Range("B16").Value2 = 45182.0
$ws.Range("C16").Value2 = 46060.0
$ws.Range("D16").Value = 'VÄRMLANDS LÄN'
$ws.Range("E16").Value = 'HAMMARÖ'
$ws.Range("F16").Value = 'Övriga Aktiebolag'
$ws.Range("G16").Value2 = 0.5
$ws.Range("H16").Value2 = 0
$ws.Range("I16").Value2 = 0
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 0
$ws.Range("L16").Value2 = 0
$ws.Range("M16").Value2 = 0
$ws.Range("N16").Value2 = 0
$ws.Range("O16").Value2 = 0
$ws.Range("P16").Value2 = 0
$ws.Range("Q16").Value2 = 0

# Row 17: A 41899-2023
$ws.Range("A17").Value = 'A 41899-2023'
$ws.Range("B17").Value2 = 45176.0
$ws.Range("C17").Value2 = 46060.0
$ws.Range("D17").Value = 'VÄRMLANDS LÄN'
$ws.Range("E17").Value = 'HAMMARÖ'
$ws.Range("F17").Value = 'Övriga Aktiebolag'
$ws.Range("G17").Value2 = 1.8
$ws.Range("H17").Value2 = 0
$ws.Range("I17").Value2 = 0
$ws.Range("J17").Value2 = 0
$ws.Range("K17").Value2 = 0
$ws.Range("L17").Value2 = 0
$ws.Range("M17").Value2 = 0
$ws.Range("N17").Value2 = 0
$ws.Range("O17").Value2 = 0
$ws.Range("P17").Value2 = 0
$ws.Range("Q17").Value2 = 0

# Row 18: A 45423-2025
$ws.Range("A18").Value = 'A 45423-2025'
$ws.Range("B18").Value2 = 45922.451377314814
$ws.Range("C18").Value2 = 46060.0
$ws.Range("D18").Value = 'VÄRMLANDS LÄN'
$ws.Range("E18").Value = 'HAMMARÖ'
$ws.Range("G18").Value2 = 1.1
$ws.Range("H18").Value2 = 0
$ws.Range("I18").Value2 = 0
$ws.Range("J18").Value2 = 0
$ws.Range("K18").Value2 = 0
$ws.Range("L18").Value2 = 0
$ws.Range("M18").Value2 = 0
$ws.Range("N18").Value2 = 0
$ws.Range("O18").Value2 = 0
$ws.Range("P18").Value2 = 0
$ws.Range("Q18").Value2 = 0

# Row 19: A 41661-2025
$ws.Range("A19").Value = 'A 41661-2025'
$ws.Range("B19").Value2 = 45902.0
$ws.Range("C19").Value2 = 46060.0
$ws.Range("D19").Value = 'VÄRMLANDS LÄN'
$ws.Range("E19").Value = 'HAMMARÖ'
$ws.Range("G19").Value2 = 5.8
$ws.Range("H19").Value2 = 0
$ws.Range("I19").Value2 = 0
$ws.Range("J19").Value2 = 0
$ws.Range("K19").Value2 = 0
$ws.Range("L19").Value2 = 0
$ws.Range("M19").Value2 = 0
$ws.Range("N19").Value2 = 0
$ws.Range("O19").Value2 = 0
$ws.Range("P19").Value2 = 0
$ws.Range("Q19").Value2 = 0

# Row 20: A 55066-2023
$ws.Range("A20").Value = 'A 55066-2023'
$ws.Range("B20").Value2 = 45237.0
$ws.Range("C20").Value2 = 46060.0
$ws.Range("D20").Value = 'VÄRMLANDS LÄN'
$ws.Range("E20").Value = 'HAMMARÖ'
$ws.Range("F20").Value = 'Kommuner'
$ws.Range("G20").Value2 = 0.9
$ws.Range("H20").Value2 = 0
$ws.Range("I20").Value2 = 0
$ws.Range("J20").Value2 = 0
$ws.Range("K20").Value2 = 0
$ws.Range("L20").Value2 = 0
$ws.Range("M20").Value2 = 0
$ws.Range("N20").Value2 = 0
$ws.Range("O20").Value2 = 0
$ws.Range("P20").Value2 = 0
$ws.Range("Q20").Value2 = 0

# Row 21: A 53369-2024
$ws.Range("A21").Value = 'A 53369-2024'
$ws.Range("B21").Value2 = 45614.0
$ws.Range("C21").Value2 = 46060.0
$ws.Range("D21").Value = 'VÄRMLANDS LÄN'
$ws.Range("E21").Value = 'HAMMARÖ'
$ws.Range("G21").Value2 = 3.4
$ws.Range("H21").Value2 = 0
$ws.Range("I21").Value2 = 0
$ws.Range("J21").Value2 = 0
$ws.Range("K21").Value2 = 0
$ws.Range("L21").Value2 = 0
$ws.Range("M21").Value2 = 0
$ws.Range("N21").Value2 = 0
$ws.Range("O21").Value2 = 0
$ws.Range("P21").Value2 = 0
$ws.Range("Q21").Value2 = 0

# Row 22: A 58383-2024
$ws.Range("A22").Value = 'A 58383-2024'
$ws.Range("B22").Value2 = 45632.0
$ws.Range("C22").Value2 = 46060.0
$ws.Range("D22").Value = 'VÄRMLANDS LÄN'
$ws.Range("E22").Value = 'HAMMARÖ'
$ws.Range("G22").Value2 = 1.6
$ws.Range("H22").Value2 = 0
$ws.Range("I22").Value2 = 0
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 0
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = 0
$ws.Range("N22").Value2 = 0
$ws.Range("O22").Value2 = 0
$ws.Range("P22").Value2 = 0
$ws.Range("Q22").Value2 = 0

# Row 23: A 2769-2023
$ws.Range("A23").Value = 'A 2769-2023'
$ws.Range("B23").Value2 = 44944.0
$ws.Range("C23").Value2 = 46060.0
$ws.Range("D23").Value = 'VÄRMLANDS LÄN'
$ws.Range("E23").Value = 'HAMMARÖ'
$ws.Range("G23").Value2 = 8.7
$ws.Range("H23").Value2 = 0
$ws.Range("I23").Value2 = 0
$ws.Range("J23").Value2 = 0
$ws.Range("K23").Value2 = 0
$ws.Range("L23").Value2 = 0
$ws.Range("M23").Value2 = 0
$ws.Range("N23").Value2 = 0
$ws.Range("O23").Value2 = 0
$ws.Range("P23").Value2 = 0
$ws.Range("Q23").Value2 = 0

# Row 24: A 42955-2023
$ws.Range("A24").Value = 'A 42955-2023'
$ws.Range("B24").Value2 = 45182.0
$ws.Range("C24").Value2 = 46060.0
$ws.Range("D24").Value = 'VÄRMLANDS LÄN'
$ws.Range("E24").Value = 'HAMMARÖ'
$ws.Range("F24").Value = 'Övriga Aktiebolag'
$ws.Range("G24").Value2 = 0.3
$ws.Range("H24").Value2 = 0
$ws.Range("I24").Value2 = 0
$ws.Range("J24").Value2 = 0
$ws.Range("K24").Value2 = 0
$ws.Range("L24").Value2 = 0
$ws.Range("M24").Value2 = 0
$ws.Range("N24").Value2 = 0
$ws.Range("O24").Value2 = 0
$ws.Range("P24").Value2 = 0
$ws.Range("Q24").Value2 = 0

# Row 25: A 42960-2023
$ws.Range("A25").Value = 'A 42960-2023'
$ws.Range("B25").Value2 = 45182.0
$ws.Range("C25").Value2 = 46060.0
$ws.Range("D25").Value = 'VÄRMLANDS LÄN'
$ws.Range("E25").Value = 'HAMMARÖ'
$ws.Range("F25").Value = 'Övriga Aktiebolag'
$ws.Range("G25").Value2 = 0.4
$ws.Range("H25").Value2 = 0
$ws.Range("I25").Value2 = 0
$ws.Range("J25").Value2 = 0
$ws.Range("K25").Value2 = 0
$ws.Range("L25").Value2 = 0
$ws.Range("M25").Value2 = 0
$ws.Range("N25").Value2 = 0
$ws.Range("O25").Value2 = 0
$ws.Range("P25").Value2 = 0
$ws.Range("Q25").Value2 = 0

# Row 26: A 55069-2023
$ws.Range("A26").Value = 'A 55069-2023'
$ws.Range("B26").Value2 = 45237.0
$ws.Range("C26").Value2 = 46060.0
$ws.Range("D26").Value = 'VÄRMLANDS LÄN'
$ws.Range("E26").Value = 'HAMMARÖ'
$ws.Range("F26").Value = 'Kommuner'
$ws.Range("G26").Value2 = 3.2
$ws.Range("H26").Value2 = 0
$ws.Range("I26").Value2 = 0
$ws.Range("J26").Value2 = 0
$ws.Range("K26").Value2 = 0
$ws.Range("L26").Value2 = 0
$ws.Range("M26").Value2 = 0
$ws.Range("N26").Value2 = 0
$ws.Range("O26").Value2 = 0
$ws.Range("P26").Value2 = 0
$ws.Range("Q26").Value2 = 0

# Row 27: A 51008-2023
$ws.Range("A27").Value = 'A 51008-2023'
$ws.Range("B27").Value2 = 45218.0
$ws.Range("C27").Value2 = 46060.0
$ws.Range("D27").Value = 'VÄRMLANDS LÄN'
$ws.Range("E27").Value = 'HAMMARÖ'
$ws.Range("F27").Value = 'Kommuner'
$ws.Range("G27").Value2 = 0.5
$ws.Range("H27").Value2 = 0
$ws.Range("I27").Value2 = 0
$ws.Range("J27").Value2 = 0
$ws.Range("K27").Value2 = 0
$ws.Range("L27").Value2 = 0
$ws.Range("M27").Value2 = 0
$ws.Range("N27").Value2 = 0
$ws.Range("O27").Value2 = 0
$ws.Range("P27").Value2 = 0
$ws.Range("Q27").Value2 = 0

# Row 28: A 41895-2023
$ws.Range("A28").Value = 'A 41895-2023'
$ws.Range("B28").Value2 = 45176.0
$ws.Range("C28").Value2 = 46060.0
$ws.Range("D28").Value = 'VÄRMLANDS LÄN'
$ws.Range("E28").Value = 'HAMMARÖ'
$ws.Range("F28").Value = 'Övriga Aktiebolag'
$ws.Range("G28").Value2 = 0.6
$ws.Range("H28").Value2 = 0
$ws.Range("I28").Value2 = 0
$ws.Range("J28").Value2 = 0
$ws.Range("K28").Value2 = 0
$ws.Range("L28").Value2 = 0
$ws.Range("M28").Value2 = 0
$ws.Range("N28").Value2 = 0
$ws.Range("O28").Value2 = 0
$ws.Range("P28").Value2 = 0
$ws.Range("Q28").Value2 = 0

# Row 29: A 13354-2022
$ws.Range("A29").Value = 'A 13354-2022'
$ws.Range("B29").Value2 = 44645.0
$ws.Range("C29").Value2 = 46060.0
$ws.Range("D29").Value = 'VÄRMLANDS LÄN'
$ws.Range("E29").Value = 'HAMMARÖ'
$ws.Range("F29").ClearContents()
$ws.Range("G29").Value2 = 1.1
$ws.Range("H29").Value2 = 0
$ws.Range("I29").Value2 = 0
$ws.Range("J29").Value2 = 0
$ws.Range("K29").Value2 = 0
$ws.Range("L29").Value2 = 0
$ws.Range("M29").Value2 = 0
$ws.Range("N29").Value2 = 0
$ws.Range("O29").Value2 = 0
$ws.Range("P29").Value2 = 0
$ws.Range("Q29").Value2 = 0

# Row 30: A 37934-2024
$ws.Range("A30").Value = 'A 37934-2024'
$ws.Range("B30").Value2 = 45544.0
$ws.Range("C30").Value2 = 46060.0
$ws.Range("D30").Value = 'VÄRMLANDS LÄN'
$ws.Range("E30").Value = 'HAMMARÖ'
$ws.Range("F30").ClearContents()
$ws.Range("G30").Value2 = 1.3
$ws.Range("H30").Value2 = 0
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 0
$ws.Range("K30").Value2 = 0
$ws.Range("L30").Value2 = 0
$ws.Range("M30").Value2 = 0
$ws.Range("N30").Value2 = 0
$ws.Range("O30").Value2 = 0
$ws.Range("P30").Value2 = 0
$ws.Range("Q30").Value2 = 0

# Row 31: A 56799-2022
$ws.Range("A31").Value = 'A 56799-2022'
$ws.Range("B31").Value2 = 44894.425625
$ws.Range("C31").Value2 = 46060.0
$ws.Range("D31").Value = 'VÄRMLANDS LÄN'
$ws.Range("E31").Value = 'HAMMARÖ'
$ws.Range("G31").Value2 = 0.4
$ws.Range("H31").Value2 = 0
$ws.Range("I31").Value2 = 0
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 0
$ws.Range("L31").Value2 = 0
$ws.Range("M31").Value2 = 0
$ws.Range("N31").Value2 = 0
$ws.Range("O31").Value2 = 0
$ws.Range("P31").Value2 = 0
$ws.Range("Q31").Value2 = 0

# Row 32: A 13356-2022
$ws.Range("A32").Value = 'A 13356-2022'
$ws.Range("B32").Value2 = 44645.0
$ws.Range("C32").Value2 = 46060.0
$ws.Range("D32").Value = 'VÄRMLANDS LÄN'
$ws.Range("E32").Value = 'HAMMARÖ'
$ws.Range("F32").ClearContents()
$ws.Range("G32").Value2 = 1.6
$ws.Range("H32").Value2 = 0
$ws.Range("I32").Value2 = 0
$ws.Range("J32").Value2 = 0
$ws.Range("K32").Value2 = 0
$ws.Range("L32").Value2 = 0
$ws.Range("M32").Value2 = 0
$ws.Range("N32").Value2 = 0
$ws.Range("O32").Value2 = 0
$ws.Range("P32").Value2 = 0
$ws.Range("Q32").Value2 = 0

# Row 33: A 42951-2023
$ws.Range("A33").Value = 'A 42951-2023'
$ws.Range("B33").Value2 = 45182.0
$ws.Range("C33").Value2 = 46060.0
$ws.Range("D33").Value = 'VÄRMLANDS LÄN'
$ws.Range("E33").Value = 'HAMMARÖ'
$ws.Range("F33").Value = 'Övriga Aktiebolag'
$ws.Range("G33").Value2 = 0.8
$ws.Range("H33").Value2 = 0
$ws.Range("I33").Value2 = 0
$ws.Range("J33").Value2 = 0
$ws.Range("K33").Value2 = 0
$ws.Range("L33").Value2 = 0
$ws.Range("M33").Value2 = 0
$ws.Range("N33").Value2 = 0
$ws.Range("O33").Value2 = 0
$ws.Range("P33").Value2 = 0
$ws.Range("Q33").Value2 = 0

# Row 34: A 57893-2023
$ws.Range("A34").Value = 'A 57893-2023'
$ws.Range("B34").Value2 = 45247.0
$ws.Range("C34").Value2 = 46060.0
$ws.Range("D34").Value = 'VÄRMLANDS LÄN'
$ws.Range("E34").Value = 'HAMMARÖ'
$ws.Range("F34").Value = 'Kommuner'
$ws.Range("G34").Value2 = 0.7
$ws.Range("H34").Value2 = 0
$ws.Range("I34").Value2 = 0
$ws.Range("J34").Value2 = 0
$ws.Range("K34").Value2 = 0
$ws.Range("L34").Value2 = 0
$ws.Range("M34").Value2 = 0
$ws.Range("N34").Value2 = 0
$ws.Range("O34").Value2 = 0
$ws.Range("P34").Value2 = 0
$ws.Range("Q34").Value2 = 0

